$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "G2"
$ws.Range("B4").Value = "Test1"
$ws.Range("C4").Value = 45860
$ws.Range("C4").Style = $ws.Range("C3").Style
$ws.Range("C4").NumberFormat = $ws.Range("C3").NumberFormat
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 0
